$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2555.22214850378
$ws.Range("C2").Value = 1743.37403972035
$ws.Range("D2").Value = 1280.88312680667
$ws.Range("F2").Value = 3681.26877844585
$ws.Range("I2").Value = 110.222148503781
$ws.Range("B3").Value = 2597.49638362041
$ws.Range("C3").Value = 1560.58769315346
$ws.Range("D3").Value = 1186.10121510877
$ws.Range("I3").Value = 1664.49638362041
$ws.Range("B4").Value = 4414.9056151499
$ws.Range("C4").Value = 3157.81603849888
$ws.Range("D4").Value = 2808.94044563282
$ws.Range("I4").Value = 906.905615149898
$ws.Range("B5").Value = 4990.37131566918
$ws.Range("I5").Value = -369.628684330823
$ws.Range("B6").Value = 5741.59408251857
$ws.Range("C6").Value = 4286.73455584564
$ws.Range("D6").Value = 3842.87151006306
$ws.Range("I6").Value = 477.594082518573
$ws.Range("B7").Value = 5984.95544467332
$ws.Range("C7").Value = 4576.30137597937
$ws.Range("I7").Value = -208.04455532668
$ws.Range("B8").Value = 5910.68186337853
$ws.Range("I8").Value = 938.68186337853
$ws.Range("B9").Value = 5129.43375580138
$ws.Range("I9").Value = 157.433755801384
$ws.Range("B10").Value = 3706.12156607305
$ws.Range("E10").Value = 5008.12892526094
$ws.Range("I10").Value = -785.878433926947
$ws.Range("B11").Value = 2748.77064436662
$ws.Range("I11").Value = -457.229355633384
$ws.Range("B12").Value = 2983.12915677683
$ws.Range("I12").Value = -627.870843223168
$ws.Range("B13").Value = 3093.39154033956
$ws.Range("I13").Value = -600.608459660443
$ws.Range("B14").Value = 2625.70802482369
$ws.Range("I14").Value = -549.291975176315
$ws.Range("B15").Value = 2633.47312104534
$ws.Range("I15").Value = -791.526878954662
$ws.Range("B16").Value = 4433.73411459297
$ws.Range("C16").Value = 2863.079709259
$ws.Range("I16").Value = -2786.26588540703
$ws.Range("B17").Value = 5026.99262560956
$ws.Range("I17").Value = -2821.00737439044
$ws.Range("B18").Value = 5779.42662112342
$ws.Range("I18").Value = -2316.57337887658
$ws.Range("B19").Value = 6004.73152874422
$ws.Range("I19").Value = -3665.26847125578
$ws.Range("B20").Value = 5941.08861073059
$ws.Range("I20").Value = -3280.91138926941
$ws.Range("B21").Value = 5245.69238773845
$ws.Range("I21").Value = -1621.30761226155
$ws.Range("B22").Value = 3915.02352487536
$ws.Range("E22").Value = 6018.59935900045
$ws.Range("I22").Value = -2016.97647512464
$ws.Range("B23").Value = 2908.29699784213
$ws.Range("I23").Value = -713.703002157868
$ws.Range("B24").Value = 3055.34389284321
$ws.Range("I24").Value = -593.656107156793
$ws.Range("B25").Value = 3149.30937091482
$ws.Range("I25").Value = -1769.69062908518
$ws.Range("B26").Value = 2696.65476749499
$ws.Range("I26").Value = -1699.34523250501
$ws.Range("B27").Value = 2709.30991602845
$ws.Range("I27").Value = -1979.69008397155
$ws.Range("B28").Value = 4476.71349127572
$ws.Range("I28").Value = -2179.28650872428
$ws.Range("B29").Value = 5082.42689782692
$ws.Range("I29").Value = -2041.57310217308
$ws.Range("B30").Value = 5815.52486278837
$ws.Range("I30").Value = -3008.47513721163
$ws.Range("B31").Value = 6053.00337511718
$ws.Range("I31").Value = -3889.99662488282
$ws.Range("B32").Value = 6004.60073919966
$ws.Range("I32").Value = -3678.39926080034
$ws.Range("B33").Value = 5333.70617534132
$ws.Range("I33").Value = -2553.29382465868
$ws.Range("B34").Value = 4134.45614157218
$ws.Range("I34").Value = -1176.54385842782
$ws.Range("B35").Value = 3127.89732682456
$ws.Range("I35").Value = 592.897326824562
$ws.Range("B36").Value = 3152.63983680369
$ws.Range("I36").Value = 583.639836803686
$ws.Range("B37").Value = 3220.4066011882
$ws.Range("I37").Value = 1400.4066011882
$ws.Range("B38").Value = 2767.72837158677
$ws.Range("I38").Value = 449.728371586774
$ws.Range("B39").Value = 2782.05728511067
$ws.Range("I39").Value = -2879.94271488933
$ws.Range("B40").Value = 4523.17101329208
$ws.Range("I40").Value = -2019.82898670792
$ws.Range("B41").Value = 5123.81189797239
$ws.Range("I41").Value = -2553.18810202761
$ws.Range("B42").Value = 5856.75102862987
$ws.Range("I42").Value = -3210.24897137013
$ws.Range("B43").Value = 6092.18601367048
$ws.Range("I43").Value = -2233.81398632952
$ws.Range("B44").Value = 6034.19512033109
$ws.Range("I44").Value = -3129.80487966891
$ws.Range("B45").Value = 5382.84815066267
$ws.Range("I45").Value = -2971.15184933733
$ws.Range("B46").Value = 4263.70931807919
$ws.Range("I46").Value = -1723.29068192081
$ws.Range("B47").Value = 3349.13548913693
$ws.Range("I47").Value = -1127.86451086307
$ws.Range("B48").Value = 3298.1434527773
$ws.Range("I48").Value = -1241.8565472227
$ws.Range("B49").Value = 3290.41384559314
$ws.Range("I49").Value = -452.58615440686
